$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "17/10/2023"
$ws.Range("B9").Value = "15:28:20"
$ws.Range("C9").Value = "15:28:25"
$ws.Range("D9").Value = "15:28:29"
$ws.Range("E9").Value = "15:28:30"
